$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (losing formatting like trailing zeros / leading zeros / fixed-point notation)
# are first forced to Text format so the literal string is preserved, matching
# the inlineStr cell type used throughout this sheet.
$textRefs = @('D4','D5','D6','D7','D8','D9','D11','D12','D13','D14','D15','D18','D19','D22','D23','D24','D25','D26','D28','D29','D30','D31','D32','D33','D34','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D45','D46','D48','D49','D51')
foreach ($ref in $textRefs) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range('D2').Value = '26.876.35'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '1.729.97'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '239.95'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').Value = '0.9996'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = '0.4838'
$ws.Range('D8').Value = '0.2597'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').Value = '0.06176'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('D10').Value = '1.728.39'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '16.05'
$ws.Range('E11').Value = '  +2.92%  '
$ws.Range('D12').Value = '0.06867'
$ws.Range('E12').Value = '  -1.87%  '
$ws.Range('D13').Value = '0.6037'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = '4.462'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').Value = '76.98'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '26.652.81'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').Value = '0.000007124'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').Value = '1.951.09'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '4.395'
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('D23').Value = '8.421'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').Value = '5.062'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').Value = '139.93'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('D26').Value = '15.19'
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('D28').Value = '106.49'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').Value = '1.376'
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('D30').Value = '3.957'
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').Value = '0.07919'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').Value = '3.664'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').Value = '0.04576'
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('D34').Value = '2.593'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').Value = '1.000'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = '0.6165'
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('D37').Value = '0.9233'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('D38').Value = '2.455'
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('D39').Value = '1.984'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').Value = '0.9992'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '5.680'
$ws.Range('E41').Value = '  +5.00%  '
$ws.Range('D42').Value = '0.01496'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('D43').Value = '99.94'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = '0.3837'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '6.783'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').Value = '0.1154'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').Value = '7.931'
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('D49').Value = '30.04'
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('E50').Value = '  +1.14%  '
$ws.Range('D51').Value = '51.26'
$ws.Range('E51').Value = '  -0.25%  '
